$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("APR-2021")
$src.Copy($null, $src)
Write-Output $wb.Worksheets.Count
foreach ($s in $wb.Worksheets) {
    Write-Output $s.Name
}
